# UndoRedoStartingStackDiagram.pptx rework
#  1. Refresh the auto "datetimeFigureOut" date placeholders (master + all
#     layouts) from 30/12/2017 to 5/3/2018.
#  2. Rename the "prevAddressBook" variable shown in the UML table to
#     "prevTaskBook" (diagram rework for the DeveloperGuide UML diagrams).

$p = $ppt.ActivePresentation

# --- 1. Date placeholders -------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "30/12/2017") {
                $tr.Text = "5/3/2018"
            }
        }
    }
}

Update-DatePlaceholders $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

# --- 2. Table text rename --------------------------------------------------
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellTr = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
                $idx = $cellTr.Text.IndexOf("prevAddressBook")
                if ($idx -ge 0) {
                    $chars = $cellTr.Characters($idx + 1, "prevAddressBook".Length)
                    $chars.Text = "prevTaskBook"
                }
            }
        }
    }
}
